$d = $word.ActiveDocument

# 1. Ativação date change
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2021", 2)

# 2. Docente responsible change
$d.Content.Find.Execute("5840560 - Marco Antonio Carvalho Pereira", $true, $false, $false, $false, $false, $true, 1, $false, "11079086 - Herlandí de Souza Andrade", 2)

# 3. Critério change
$d.Content.Find.Execute("Média aritmética de duas provas teóricas.", $true, $false, $false, $false, $false, $true, 1, $false, "Média Aritmética das atividades avaliativas realizadas.", 2)

# 4. Norma de recuperação - remove trailing period
$d.Content.Find.Execute("Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação", 2)

# 5. Add new requisito line after LOQ4205 line
$d.Content.Find.Execute("LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`rLOQ4240 -  Administração e Organização II  (Requisito fraco)", 2)

$word.Documents
